# Update existing annual GDP observations (col B) with revised FRED figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 102596.12699999999
$ws.Range("B13").Value = 105867.63499999999
$ws.Range("B14").Value = 111060.939
$ws.Range("B15").Value = 116501.16
$ws.Range("B16").Value = 122524.352
$ws.Range("B17").Value = 126569.159
$ws.Range("B18").Value = 131097.4
$ws.Range("B19").Value = 136423.61499999999
$ws.Range("B20").Value = 135802.17300000001
$ws.Range("B21").Value = 139383.24400000001
$ws.Range("B22").Value = 142235.693
$ws.Range("B23").Value = 147315.27100000001
$ws.Range("B24").Value = 150924.34599999999
$ws.Range("B25").Value = 156102.29699999999
$ws.Range("B26").Value = 159904.74
$ws.Range("B27").Value = 160781.04199999999
$ws.Range("B28").Value = 163129.43700000001
$ws.Range("B29").Value = 170039.45
$ws.Range("B30").Value = 175382.47

# Append the new 2020-01-01 observation as row 31, matching the formatting
# (date style) used by the existing date column and copying the 0.000
# number style used by the existing GDP column.
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("B30").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A31").Value = 43831
$ws.Range("B31").Value = 171493.44399999999

# Mirror the workbook's on-save cursor position: whole columns A:B selected
# with the active cell resting on A14.
$ws.Range("A:B").Select()
